$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.600934028625488
$ws.Range("B1").Value = 3.664695024490356
$ws.Range("C1").Value = 3.207427024841309
$ws.Range("D1").Value = 3.976921796798706
$ws.Range("E1").Value = 5.030847072601318
